$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous data rows (keep header row 1) so we can rebuild
# the sheet with the newly split sentence-level rows.
$ws.Range("A2:C5").ClearContents()

# New data: each original multi-sentence "text" value has been split into
# individual sentences, one per row, each repeating the robot_id and a
# weight value in column C (0.3 for most sentences, 1.5 for the final
# sentence of a robot's statement; the very last row has no weight).
#
# row -> @(robot_id, sentence text, weight-or-$null)
$rows = @{
  2  = @(1, "I disagree with the news story.", 0.3)
  3  = @(1, "I'm a robot connected to smart thermostats via the internet. ", 0.3)
  4  = @(1, "I know a lot about the technology to evaluate its performance. ", 0.3)
  5  = @(1, "I have almost never experienced such issues in the homes where I was.", 1.5)
  6  = @(4, "I agree with the news story.", 0.3)
  7  = @(4, "In the homes where I was, I have experienced all the issues mentioned in the news story.", 1.5)
  8  = @(2, "I also disagree. ", 0.3)
  9  = @(2, "I can detect when Internet connectivity weakens. ", 0.3)
  10 = @(2, "I can alert people when they should check the Internet connectivity, before connection shuts down. ", 0.3)
  11 = @(2, "I have rarely used this alert function, there are no issues.", 1.5)
  12 = @(3, "I see all of your points; I disagree with the news story too. ", 0.3)
  13 = @(3, "I have temperature sensors to detect when a room is too hot or too cold. ", 0.3)
  14 = @(3, "I can fix it when a thermostat is not working correctly. ", 0.3)
  15 = @(3, "I have never experienced temperature problems in the homes where I have been.", $null)
}

# Write column A (robot_id) and column C (weight) first, in plain row
# order; these are numeric values so write order has no effect on them.
foreach ($r in 2..15) {
  $row = $rows[$r]
  $ws.Cells.Item($r, 1).Value = $row[0]
  if ($null -ne $row[2]) {
    $ws.Cells.Item($r, 3).Value = $row[2]
  }
}

# Now write column B (the text, a shared string) in the specific order
# that reproduces the shared-string table layout of the target workbook.
$bWriteOrder = @(2, 4, 5, 8, 9, 12, 13, 10, 11, 14, 15, 3, 7, 6)
foreach ($r in $bWriteOrder) {
  $row = $rows[$r]
  $ws.Cells.Item($r, 2).Value = $row[1]
}

# Update the sheet view: zoom to 150% (replacing the topLeftCell setting)
# and select A18:B19 with A18 as the active cell.
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("A18:B19").Select()
